# "level up exp config" -- add a new "LevelupExp" sheet with a
# level -> required-exp table, right after the existing HeroProto sheet,
# make it the active sheet/tab, move HeroProto's own selection to A2,
# and drop the now-stale "1 客户端专用 / 2 服务器专用 / ..." author note
# that lived on HeroProto!A2.

$wb = $excel.ActiveWorkbook
$hero = $wb.Worksheets.Item(1)

# --- add the new worksheet after HeroProto -------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "LevelupExp"

# --- header rows (mirrors the 4-row HeroProto schema-header convention) --
$ws.Range("A1").Value = "等级"
$ws.Range("B1").Value = "经验"
$ws.Range("A2").Value = "jl"
$ws.Range("B2").Value = "jl"
$ws.Range("A3").Value = "level"
$ws.Range("B3").Value = "exp"
$ws.Range("A4").Value = "number"
$ws.Range("B4").Value = "number"

# --- match HeroProto's header formatting (fonts/styles) on the new sheet -
$fmt = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

$hero.Range("B1").Copy()
$ws.Range("A1").PasteSpecial($fmt)
$hero.Range("G1").Copy()
$ws.Range("B1").PasteSpecial($fmt)

$hero.Range("B2").Copy()
$ws.Range("A2").PasteSpecial($fmt)
$ws.Range("B2").PasteSpecial($fmt)

$hero.Range("B3").Copy()
$ws.Range("A3").PasteSpecial($fmt)
$ws.Range("B3").PasteSpecial($fmt)

$hero.Range("B4").Copy()
$ws.Range("A4").PasteSpecial($fmt)
$hero.Range("A4").Copy()
$ws.Range("B4").PasteSpecial($fmt)

$excel.CutCopyMode = $false

# --- level/exp data rows 5-29 (levels 1-25) -------------------------------
$level = 1
for ($row = 5; $row -le 29; $row++) {
    $ws.Cells.Item($row, 1).Value = $level
    $ws.Cells.Item($row, 2).Value = 1000 * $level + 1000
    $level++
}

# --- drop the stale maintenance note on HeroProto!A2 ----------------------
$hero.Range("A2").Comment.Delete()

# --- final selections: HeroProto cursor moves to A2, LevelupExp -> K20,
#     and LevelupExp becomes the active/visible tab -----------------------
$hero.Range("A2").Select()
$ws.Range("K20").Select()
